$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.275.94'
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").Value = '1.865.76'
$ws.Range("E3").Value = '  +0.44%  '

$ws.Range("D4").Value = '''1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''236.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.71%  '

$ws.Range("D6").Value = '''1.0000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = '''0.4680'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.42%  '

$ws.Range("D8").Value = '''0.2865'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.14%  '

$ws.Range("D9").Value = '''0.06545'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.35%  '

$ws.Range("D10").Value = '''22.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +14.90%  '

$ws.Range("D11").Value = '''0.07909'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.25%  '

$ws.Range("D12").Value = '''97.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.77%  '

$ws.Range("D13").Value = '1.870.60'
$ws.Range("E13").Value = '  +0.67%  '

$ws.Range("D14").Value = '''5.179'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.43%  '

$ws.Range("D15").Value = '''0.6826'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.15%  '

$ws.Range("D16").Value = '''279.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("D17").Value = '30.271.72'
$ws.Range("E17").Value = '  +0.53%  '

$ws.Range("D18").Value = '''13.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.27%  '

$ws.Range("D19").Value = '''0.9999'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '

$ws.Range("D20").Value = '''5.411'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.01%  '

$ws.Range("D21").Value = '''0.000007341'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.86%  '

$ws.Range("D22").Value = '2.113.14'
$ws.Range("E22").Value = '  +0.56%  '

$ws.Range("D23").Value = '''1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").Value = '''6.191'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.45%  '

$ws.Range("D25").Value = '''168.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.25%  '

$ws.Range("D26").Value = '''9.275'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").Value = '''19.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.09%  '

$ws.Range("D28").Value = '''1.943'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.61%  '

$ws.Range("D29").Value = '''1.382'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.86%  '

$ws.Range("D30").Value = '''0.09851'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.42%  '

$ws.Range("D31").Value = '''4.402'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("D32").Value = '''1.482'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.32%  '

$ws.Range("D33").Value = '''4.074'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D34").Value = '''0.04754'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.50%  '

$ws.Range("D35").Value = '''1.140'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.80%  '

$ws.Range("D36").Value = '''0.7118'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.96%  '

$ws.Range("D37").Value = '''2.708'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.45%  '

$ws.Range("D38").Value = '''0.01881'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.44%  '

$ws.Range("D39").Value = '''2.614'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.05%  '

$ws.Range("D40").Value = '''77.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.44%  '

$ws.Range("D41").Value = '''6.308'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.82%  '

$ws.Range("D42").Value = '''1.962'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.85%  '

$ws.Range("D43").Value = '''0.8529'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("D44").Value = '''0.4194'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.62%  '

$ws.Range("D45").Value = '''0.9994'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").Value = '''103.55'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("D47").Value = '''963.73'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.08%  '

$ws.Range("D48").Value = '''7.236'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.11%  '

$ws.Range("D49").Value = '''9.293'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.26%  '

$ws.Range("D50").Value = '''34.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.81%  '

$ws.Range("D51").Value = '''0.05644'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.35%  '
